# "c++ hashing, java trees"
# Applies to the "2024 Fall" sheet (Mixed_Review.xlsx, java library):
#  1. Update the notes for "198. House Robber" (row 4) to mention Fibonacci.
#  2. Append two new rows to the review table for:
#       226. Invert Binary Tree
#       104. Maximum Depth of Binary Tree
#  3. Grow the table / dimension to match, and add the new hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024 Fall")

# 1) Tweak the existing House Robber note (row 4, column D).
$ws.Range("D4").Value = "Recursive relation argmax(dfs(i-2) + curr, dfs(i-1)). Fibonacci is Optimal."

# 2) New row 7: 226. Invert Binary Tree
$ws.Range("E7").Value = "https://leetcode.com/problems/invert-binary-tree/solutions/62707/straightforward-dfs-recursive-iterative-bfs-solutions/"
$ws.Hyperlinks.Add($ws.Range("E7"), "https://leetcode.com/problems/invert-binary-tree/solutions/62707/straightforward-dfs-recursive-iterative-bfs-solutions/") | Out-Null
$ws.Range("E7").Style = "Hyperlink"
$ws.Range("D7").Value = "Recursive DFS > Stack > Iterative BFS (Level Order Traversal)"
$ws.Range("A7").Value = "226. Invert Binary Tree"
$ws.Range("B7").Value = "Easy"
$ws.Range("B7").Interior.Color = 5287936
$ws.Range("C7").Value = "Trees"

# 3) New row 8: 104. Maximum Depth of Binary Tree
$ws.Range("A8").Value = "104. Maximum Depth of Binary Tree"
$ws.Range("B8").Value = "Easy"
$ws.Range("B8").Interior.Color = 5287936
$ws.Range("C8").Value = "Trees"
$ws.Range("E8").Value = "https://leetcode.com/problems/maximum-depth-of-binary-tree/solutions/1770060/c-recursive-dfs-example-dry-run-well-explained/ "
$ws.Hyperlinks.Add($ws.Range("E8"), "https://leetcode.com/problems/maximum-depth-of-binary-tree/solutions/1770060/c-recursive-dfs-example-dry-run-well-explained/ ") | Out-Null
$ws.Range("E8").Style = "Hyperlink"
$ws.Range("D8").Value = "Recursive DFS of both sides, and recursive relation of max(maxL, maxR)+1;"

# Grow the table to include the two new rows.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E8"))

# Match the saved selection shown in the diff.
$ws.Range("E13").Select() | Out-Null
